$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Insert the GitHub repo link into the blank "NoSpacing" paragraph
#    that currently sits between the assignment title and the
#    screenshot image, then leave a fresh blank "NoSpacing" paragraph
#    behind it (so the image paragraph keeps its own line).
# -------------------------------------------------------------------

# Locate the existing blank paragraph right before the picture
# paragraph (4th paragraph in the document body).
$blankPara = $d.Paragraphs.Item(4)
$blankRange = $blankPara.Range
$blankRange.Collapse(0)

# Split off a brand new blank paragraph first (inheriting the
# NoSpacing style) so the original paragraph keeps its identity and
# becomes the one that carries the hyperlink.
$blankRange.InsertParagraphAfter()

# Re-fetch paragraph 4 (still the original, now-hyperlinked paragraph)
# and append the hyperlink run into it.
$linkPara = $d.Paragraphs.Item(4)
$linkRange = $linkPara.Range
$linkRange.Collapse(0)
$d.Hyperlinks.Add($linkRange, "https://github.com/trueworthy/csd-340/tree/main/module-1")

# -------------------------------------------------------------------
# 2) Register the "Unresolved Mention" latent character style that
#    Word's collaborative-comment tooling adds alongside the edit.
# -------------------------------------------------------------------
$mentionStyle = $d.Styles.Add("Unresolved Mention", 2)
$mentionStyle.BaseStyle = "DefaultParagraphFont"
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionStyle.Font.Color = 6053472
$mentionStyle.Font.Shading.BackgroundPatternColor = 14540769

Write-Output "Applied hyperlink + Unresolved Mention style."
